# Apply the "trimestral" invoice sheet update:
#  - Header tweaks (H2, I2, J2)
#  - Replace test/placeholder invoice rows (3-13) with real invoice data
#  - Delete row 14 (only 11 real invoices remain, rows 3-13)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row tweaks ---
$ws.Range("H2").Value = "IVA %"
$ws.Range("I2").Value = "Descuento % "
$ws.Range("J2").Value = "Total facturado €(Impuestos y descuentos incluidos"

# --- Delete row 14 first so everything below shifts up cleanly ---
$ws.Rows.Item(14).Delete()

# Date-like text columns (B, E) must stay text, not get auto-converted to
# Excel date serials, so force a text number format before assignment.
$dateCols = "B3","E3","B4","E4","B5","E5","B6","E6","B7","E7","B8","E8","B9","E9","B10","E10","B11","E11","B12","E12","B13","E13"
foreach ($addr in $dateCols) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Row 3: A200001 / Clientes contado ---
$ws.Range("A3").Value = "A200001"
$ws.Range("B3").Value = "2020-09-06"
$ws.Range("C3").Value = "Clientes contado "
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "2020-09-06"
$ws.Range("F3").Value = "Sesión de Fisioterapia"
$ws.Range("G3").Value = 32
$ws.Range("H3").Value = 4
$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 31.62

# --- Row 4: A200002 / Luis Zurita Herrera ---
$ws.Range("A4").Value = "A200002"
$ws.Range("B4").Value = "2020-09-06"
$ws.Range("C4").Value = "Luis Zurita Herrera"
$ws.Range("D4").Value = "74666101M"
$ws.Range("E4").Value = "2020-09-06"
$ws.Range("F4").Value = "Sesión de Fisioterapia"
$ws.Range("G4").Value = 64
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 64

# --- Row 5: A200003 / Luis Zurita Herrera ---
$ws.Range("A5").Value = "A200003"
$ws.Range("B5").Value = "2020-09-06"
$ws.Range("C5").Value = "Luis Zurita Herrera"
$ws.Range("D5").Value = "74666101M"
$ws.Range("E5").Value = "2020-09-06"
$ws.Range("F5").Value = "Sesión de Fisioterapia"
$ws.Range("G5").Value = 32
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 15
$ws.Range("J5").Value = 27.2

# --- Row 6: A200004 / Clientes contado ---
$ws.Range("A6").Value = "A200004"
$ws.Range("B6").Value = "2020-09-06"
$ws.Range("C6").Value = "Clientes contado "
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = "2020-09-06"
$ws.Range("F6").Value = "Sesión de Fisioterapia"
$ws.Range("G6").Value = 32
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 32

# --- Row 7: A200005 / Clientes contado ---
$ws.Range("A7").Value = "A200005"
$ws.Range("B7").Value = "2020-09-06"
$ws.Range("C7").Value = "Clientes contado "
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = "2020-09-06"
$ws.Range("F7").Value = "Sesión de Fisioterapia"
$ws.Range("G7").Value = 64
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 64

# --- Row 8: A200006 / Clientes contado ---
$ws.Range("A8").Value = "A200006"
$ws.Range("B8").Value = "2020-09-06"
$ws.Range("C8").Value = "Clientes contado "
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = "2020-09-06"
$ws.Range("F8").Value = "Sesión de Fisioterapia"
$ws.Range("G8").Value = 32
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 32

# --- Row 9: A200007 / Clientes contado / Acupuntura ---
$ws.Range("A9").Value = "A200007"
$ws.Range("B9").Value = "2020-09-06"
$ws.Range("C9").Value = "Clientes contado "
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = "2020-09-06"
$ws.Range("F9").Value = "Sesión de Acupuntura"
$ws.Range("G9").Value = 81.59
$ws.Range("H9").Value = 21
$ws.Range("I9").Value = 23
$ws.Range("J9").Value = 76.02

# --- Row 10: A200008 / Clientes contado (amounts unchanged) ---
$ws.Range("A10").Value = "A200008"
$ws.Range("B10").Value = "2020-09-06"
$ws.Range("C10").Value = "Clientes contado "
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = "2020-09-06"
$ws.Range("F10").Value = "ffdfasdfdasfdsaf"
$ws.Range("G10").Value = 32
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 32

# --- Row 11: A200009 / Clientes contado / acu ---
$ws.Range("A11").Value = "A200009"
$ws.Range("B11").Value = "2020-09-06"
$ws.Range("C11").Value = "Clientes contado "
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = "2020-09-06"
$ws.Range("F11").Value = "acu"
$ws.Range("G11").Value = 50
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 50

# --- Row 12: A200010 / Luis Zurita Herrera / acu ---
$ws.Range("A12").Value = "A200010"
$ws.Range("B12").Value = "2020-09-06"
$ws.Range("C12").Value = "Luis Zurita Herrera"
$ws.Range("D12").Value = "74666101M"
$ws.Range("E12").Value = "2020-09-06"
$ws.Range("F12").Value = "acu"
$ws.Range("G12").Value = 100
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 12
$ws.Range("J12").Value = 88

# --- Row 13: A200011 / Clientes contado ---
$ws.Range("A13").Value = "A200011"
$ws.Range("B13").Value = "2020-09-06"
$ws.Range("C13").Value = "Clientes contado "
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = "2020-09-06"
$ws.Range("F13").Value = "Sesión de Fisioterapia"
$ws.Range("G13").Value = 32
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 32
